$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (A238) down
# through the new rows A239:A244 so the new date cells keep the same
# bold/centered/bordered date style used throughout column A.
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)

# New daily records continuing the series through 2021-05-02.
$ws.Range("A239").Value = 44313
$ws.Range("B239").Value = 0
$ws.Range("C239").Value = 0
$ws.Range("D239").Value = 0

$ws.Range("A240").Value = 44314
$ws.Range("B240").Value = 0
$ws.Range("C240").Value = 0
$ws.Range("D240").Value = 0

$ws.Range("A241").Value = 44315
$ws.Range("B241").Value = 0
$ws.Range("C241").Value = 0
$ws.Range("D241").Value = 0

$ws.Range("A242").Value = 44316
$ws.Range("B242").Value = 1
$ws.Range("C242").Value = 1
$ws.Range("D242").Value = 37.46721618583739

$ws.Range("A243").Value = 44317
$ws.Range("B243").Value = 0
$ws.Range("C243").Value = 1
$ws.Range("D243").Value = 37.46721618583739

$ws.Range("A244").Value = 44318
$ws.Range("B244").Value = 1
$ws.Range("C244").Value = 2
$ws.Range("D244").Value = 74.93443237167479
